# LoginStatus.xlsx — add a new "status" column (L) that mirrors the
# existing column K (header "status" + PASS/FAIL rows), matching the
# thread-local-driver test run that was appended to the login-status
# report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column K (11) already holds the "status" header in row 1 and PASS/FAIL
# values below it for every test run so far. Mirror that whole column
# into the new column L (12) for the additional run.
$srcCol = 11
$dstCol = 12

$lastRow = $ws.Cells.Item($ws.Rows.Count, $srcCol).End(-4162).Row
if ($lastRow -lt 6) {
    $lastRow = 6
}

for ($r = 1; $r -le $lastRow; $r++) {
    $src = $ws.Cells.Item($r, $srcCol)
    $dst = $ws.Cells.Item($r, $dstCol)
    $dst.Value = $src.Value2
}

# Re-apply the same header fill (the grey "status" header shading used by
# every other column) to the new header cell L1.
$ws.Cells.Item(1, $dstCol).Interior.ColorIndex = $ws.Cells.Item(1, $srcCol).Interior.ColorIndex

# Match column K's width for the new column.
$ws.Columns($dstCol).ColumnWidth = $ws.Columns($srcCol).ColumnWidth

Write-Host "Added column L mirroring column K (status/PASS/FAIL)."
